# "change excel2json to many" - monitor_xlsx/20260129.xlsx update
#
# Updates the foreign-investor 20-day average/sum figures (and the
# P/C-ratio roll-up) on the overview + detail sheets, and refreshes the
# per-stock chip-monitor sheet: new MA20-deviation values, a few
# recomputed price/volume figures, and drops the now-unavailable
# broker-breadth / short-covering / VWAP columns (P:U) in favour of a
# plain "N/A" data-source marker.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item(1)   # 總覽
$wsDetail   = $wb.Worksheets.Item(2)   # 詳細數據
$wsStocks   = $wb.Worksheets.Item(3)   # 個股籌碼

# ---------------------------------------------------------------------
# 總覽 (Overview)
# ---------------------------------------------------------------------
$wsOverview.Range("G9").Value = "52.03億"
$wsOverview.Range("H9").Value = "1040.59億"
# Leading apostrophe = Excel's own "force text" quote-prefix so the
# "...%" strings stay literal text instead of being auto-parsed into a
# percentage number (matches the source cells, which are plain text).
$wsOverview.Range("C11").Value = "'158.39%"
$wsOverview.Range("E11").Value = "'170.5%"

# ---------------------------------------------------------------------
# 詳細數據 (Detail)
# ---------------------------------------------------------------------
$wsDetail.Range("B6").Value = "52.03億"
$wsDetail.Range("B7").Value = "1040.59億"
$wsDetail.Range("B21").Value = "'170.5%"

# ---------------------------------------------------------------------
# 個股籌碼 (Per-stock chips)
# ---------------------------------------------------------------------

# Row 4 - 0050
$wsStocks.Range("N4").Value = -0.79
$wsStocks.Range("P4:U4").Clear()
$wsStocks.Range("V4").Value = "N/A"

# Row 5 - 00708L (no N column)
$wsStocks.Range("P5:U5").Clear()
$wsStocks.Range("V5").Value = "N/A"

# Row 6 - 1519
$wsStocks.Range("N6").Value = 21.14
$wsStocks.Range("P6:U6").Clear()
$wsStocks.Range("V6").Value = "N/A"

# Row 7 - 1605
$wsStocks.Range("N7").Value = 23.04
$wsStocks.Range("P7:U7").Clear()
$wsStocks.Range("V7").Value = "N/A"

# Row 8 - 2308
$wsStocks.Range("N8").Value = 21.82
$wsStocks.Range("P8:U8").Clear()
$wsStocks.Range("V8").Value = "N/A"

# Row 9 - 2330
$wsStocks.Range("N9").Value = 23.09
$wsStocks.Range("P9:U9").Clear()
$wsStocks.Range("V9").Value = "N/A"

# Row 10 - 2344
$wsStocks.Range("N10").Value = 37.67
$wsStocks.Range("P10:U10").Clear()
$wsStocks.Range("V10").Value = "N/A"

# Row 11 - 2383
$wsStocks.Range("N11").Value = 24.4
$wsStocks.Range("P11:U11").Clear()
$wsStocks.Range("V11").Value = "N/A"

# Row 12 - 3661
$wsStocks.Range("N12").Value = 19.49
$wsStocks.Range("P12:U12").Clear()
$wsStocks.Range("V12").Value = "N/A"

# Row 13 - 4958
$wsStocks.Range("G13").Value = 4701
$wsStocks.Range("I13").Value = -1119
$wsStocks.Range("N13").Value = 19.61
$wsStocks.Range("O13").Value = "中性"
$wsStocks.Range("P13:U13").Clear()
$wsStocks.Range("V13").Value = "N/A"

# Row 14 - 6442
$wsStocks.Range("N14").Value = 55.09
$wsStocks.Range("P14:U14").Clear()
$wsStocks.Range("V14").Value = "N/A"

# Row 15 - 3081
$wsStocks.Range("C15").Value = 1000
$wsStocks.Range("D15").Value = -2.44
$wsStocks.Range("E15").Value = 905
$wsStocks.Range("P15:U15").Clear()
$wsStocks.Range("V15").Value = "N/A"

# Row 16 - 3260 (D16 also flips from the "up/green" style to the
# "down/red" style used elsewhere in the sheet, matching its new
# negative value)
$wsStocks.Range("C16").Value = 306.5
$wsStocks.Range("D4").Copy()
$wsStocks.Range("D16").PasteSpecial(-4122)
$wsStocks.Range("D16").Value = -8.1
$wsStocks.Range("E16").Value = 26486
$wsStocks.Range("P16:U16").Clear()
$wsStocks.Range("V16").Value = "N/A"

# Row 17 - 3265
$wsStocks.Range("C17").Value = 132.5
$wsStocks.Range("D17").Value = -2.57
$wsStocks.Range("E17").Value = 1591
$wsStocks.Range("P17:U17").Clear()
$wsStocks.Range("V17").Value = "N/A"

# Row 18 - 4979
$wsStocks.Range("C18").Value = 318.5
$wsStocks.Range("D18").Value = -0.16
$wsStocks.Range("E18").Value = 20497
$wsStocks.Range("P18:U18").Clear()
$wsStocks.Range("V18").Value = "N/A"

# Row 19 - 3189
$wsStocks.Range("N19").Value = 23.88
$wsStocks.Range("P19:U19").Clear()
$wsStocks.Range("V19").Value = "N/A"
